# Daily attendance processing - 2025-12-26 14:06:32
# Reorders the "Recorded By" (column G) entries so that "System" is listed
# first among the comma-separated recorders, for every row where the list
# of recorders ends with "System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; Value='System, system, backup@backdoor.com'},
    @{Row=3; Value='System, dnasr281@gmail.com'},
    @{Row=4; Value='System, backup@backdoor.com'},
    @{Row=5; Value='System, backup@backdoor.com'},
    @{Row=6; Value='System, dnasr281@gmail.com'},
    @{Row=7; Value='System, admin@admin.com'},
    @{Row=8; Value='System, backup@backdoor.com'},
    @{Row=10; Value='System, dnasr281@gmail.com'},
    @{Row=11; Value='System, dnasr281@gmail.com'},
    @{Row=12; Value='System, dnasr281@gmail.com'},
    @{Row=13; Value='System, dnasr281@gmail.com'},
    @{Row=14; Value='System, dnasr281@gmail.com'},
    @{Row=15; Value='System, dnasr281@gmail.com'},
    @{Row=17; Value='System, dnasr281@gmail.com'},
    @{Row=18; Value='System, dnasr281@gmail.com'},
    @{Row=19; Value='System, dnasr281@gmail.com'},
    @{Row=20; Value='System, dnasr281@gmail.com'},
    @{Row=21; Value='System, dnasr281@gmail.com'},
    @{Row=22; Value='System, dnasr281@gmail.com'},
    @{Row=24; Value='System, dnasr281@gmail.com'},
    @{Row=26; Value='System, dnasr281@gmail.com'},
    @{Row=28; Value='System, system, backup@backdoor.com'},
    @{Row=29; Value='System, dnasr281@gmail.com'},
    @{Row=30; Value='System, backup@backdoor.com'},
    @{Row=31; Value='System, backup@backdoor.com'},
    @{Row=32; Value='System, dnasr281@gmail.com'},
    @{Row=33; Value='System, admin@admin.com'},
    @{Row=34; Value='System, backup@backdoor.com'},
    @{Row=36; Value='System, dnasr281@gmail.com'},
    @{Row=37; Value='System, dnasr281@gmail.com'},
    @{Row=38; Value='System, dnasr281@gmail.com'},
    @{Row=39; Value='System, dnasr281@gmail.com'},
    @{Row=40; Value='System, dnasr281@gmail.com'},
    @{Row=41; Value='System, dnasr281@gmail.com'},
    @{Row=43; Value='System, dnasr281@gmail.com'},
    @{Row=44; Value='System, dnasr281@gmail.com'},
    @{Row=45; Value='System, dnasr281@gmail.com'},
    @{Row=46; Value='System, dnasr281@gmail.com'},
    @{Row=47; Value='System, dnasr281@gmail.com'},
    @{Row=48; Value='System, dnasr281@gmail.com'},
    @{Row=50; Value='System, dnasr281@gmail.com'},
    @{Row=52; Value='System, dnasr281@gmail.com'},
    @{Row=54; Value='System, system, backup@backdoor.com'},
    @{Row=55; Value='System, dnasr281@gmail.com'},
    @{Row=56; Value='System, backup@backdoor.com'},
    @{Row=57; Value='System, backup@backdoor.com'},
    @{Row=58; Value='System, dnasr281@gmail.com'},
    @{Row=59; Value='System, admin@admin.com'},
    @{Row=60; Value='System, backup@backdoor.com'},
    @{Row=62; Value='System, dnasr281@gmail.com'},
    @{Row=63; Value='System, dnasr281@gmail.com'},
    @{Row=64; Value='System, dnasr281@gmail.com'},
    @{Row=65; Value='System, dnasr281@gmail.com'},
    @{Row=66; Value='System, dnasr281@gmail.com'},
    @{Row=67; Value='System, dnasr281@gmail.com'},
    @{Row=69; Value='System, dnasr281@gmail.com'},
    @{Row=70; Value='System, dnasr281@gmail.com'},
    @{Row=71; Value='System, dnasr281@gmail.com'},
    @{Row=72; Value='System, dnasr281@gmail.com'},
    @{Row=73; Value='System, dnasr281@gmail.com'},
    @{Row=74; Value='System, dnasr281@gmail.com'},
    @{Row=76; Value='System, dnasr281@gmail.com'},
    @{Row=78; Value='System, dnasr281@gmail.com'},
    @{Row=80; Value='System, backup@backdoor.com'},
    @{Row=81; Value='System, backup@backdoor.com'},
    @{Row=82; Value='System, backup@backdoor.com'},
    @{Row=83; Value='System, dnasr281@gmail.com'},
    @{Row=84; Value='System, dnasr281@gmail.com'},
    @{Row=85; Value='System, dnasr281@gmail.com'},
    @{Row=86; Value='System, dnasr281@gmail.com'},
    @{Row=90; Value='System, dnasr281@gmail.com'},
    @{Row=92; Value='System, dnasr281@gmail.com'},
    @{Row=93; Value='System, dnasr281@gmail.com'},
    @{Row=94; Value='System, dnasr281@gmail.com'},
    @{Row=96; Value='System, dnasr281@gmail.com'},
    @{Row=99; Value='System, dnasr281@gmail.com'},
    @{Row=101; Value='System, dnasr281@gmail.com'},
    @{Row=106; Value='System, backup@backdoor.com'},
    @{Row=107; Value='System, backup@backdoor.com'},
    @{Row=108; Value='System, backup@backdoor.com'},
    @{Row=109; Value='System, dnasr281@gmail.com'},
    @{Row=110; Value='System, dnasr281@gmail.com'},
    @{Row=111; Value='System, dnasr281@gmail.com'},
    @{Row=112; Value='System, dnasr281@gmail.com'},
    @{Row=116; Value='System, dnasr281@gmail.com'},
    @{Row=118; Value='System, dnasr281@gmail.com'},
    @{Row=119; Value='System, dnasr281@gmail.com'},
    @{Row=120; Value='System, dnasr281@gmail.com'},
    @{Row=122; Value='System, dnasr281@gmail.com'},
    @{Row=125; Value='System, dnasr281@gmail.com'},
    @{Row=127; Value='System, dnasr281@gmail.com'},
    @{Row=132; Value='System, backup@backdoor.com'},
    @{Row=133; Value='System, backup@backdoor.com'},
    @{Row=134; Value='System, backup@backdoor.com'},
    @{Row=135; Value='System, dnasr281@gmail.com'},
    @{Row=136; Value='System, dnasr281@gmail.com'},
    @{Row=137; Value='System, dnasr281@gmail.com'},
    @{Row=138; Value='System, dnasr281@gmail.com'},
    @{Row=142; Value='System, dnasr281@gmail.com'},
    @{Row=144; Value='System, dnasr281@gmail.com'},
    @{Row=145; Value='System, dnasr281@gmail.com'},
    @{Row=146; Value='System, dnasr281@gmail.com'},
    @{Row=148; Value='System, dnasr281@gmail.com'},
    @{Row=151; Value='System, dnasr281@gmail.com'},
    @{Row=153; Value='System, dnasr281@gmail.com'}
)

foreach ($update in $updates) {
    $ws.Cells.Item($update.Row, 7).Value = $update.Value
}
